# update after 2nd teaching
# 1) Slide 1 ("Subtitle 2" placeholder): change the second bullet line.
$p = $ppt.ActivePresentation

$s1 = $p.Slides.Item(1)
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $shp = $s1.Shapes.Item($i)
    if ($shp.Name -eq "Subtitle 2" -and $shp.HasTextFrame) {
        $tr = $shp.TextFrame.TextRange
        $tr.Replace("use login/logout for question answering", "class view")
    }
}

# 2) Refresh the cached "datetimeFigureOut" footer date field (11/16/2021 -> 12/23/2021)
#    on the slide master and on every slide layout.
$newDate = "12/23/2021"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $shp = $master.Shapes.Item($i)
    if ($shp.Type -eq 14 -and $shp.HasTextFrame) {
        if ($shp.PlaceholderFormat.Type -eq 16) {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $shp = $layout.Shapes.Item($i)
        if ($shp.Type -eq 14 -and $shp.HasTextFrame) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}
